# Update the public EPEX Spot prices workbook with the latest day of data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": append a new day column (AU) with header "30-jul"
# and 24 hourly price values.
# ---------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous header cell (bold, centered, bordered)
# onto the new header cell, then set its text.
$wsSpot.Range("AT1").Copy()
$wsSpot.Range("AU1").PasteSpecial(-4122)
$wsSpot.Range("AU1").Value = "30-jul"

$spotValues = @(
    54.44,
    50.92,
    43.75,
    38.48,
    32.55,
    38.48,
    45.51,
    72.95999999999999,
    57.38,
    46.74,
    33.59,
    39.88,
    24.79,
    17.69,
    10.7,
    24.21,
    19.89,
    31.06,
    40.75,
    69.55,
    77.98999999999999,
    101.49,
    104.36,
    92.18000000000001
)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 47).Value = $spotValues[$i]
}

# ---------------------------------------------------------------------
# Sheet "Gaz": append the next day's closing price.
# The date column is stored as plain text, so force a text format while
# assigning the value (otherwise it gets auto-parsed into a date serial),
# then clear the temporary formatting back to the sheet's default style.
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A44").NumberFormat = "@"
$wsGaz.Range("A44").Value = "2025-07-28"
$wsGaz.Range("A44").ClearFormats()
$wsGaz.Range("B44").Value = 32.5

# ---------------------------------------------------------------------
# Sheet "CO2": append the next day's closing price.
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A44").NumberFormat = "@"
$wsCo2.Range("A44").Value = "2025-07-28"
$wsCo2.Range("A44").ClearFormats()
$wsCo2.Range("B44").Value = 69.5
